$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-04 Friday" "2024-10-05 Saturday"

Replace-Text "77×26=" "75×37="
Replace-Text "20×56=" "60×81="
Replace-Text "87×26=" "38×66="
Replace-Text "27×38=" "61×95="
Replace-Text "97×53=" "46×40="

Replace-Text "77×81=" "69×17="
Replace-Text "28×66=" "68×13="
Replace-Text "51×45=" "55×51="
Replace-Text "24×12=" "39×61="
Replace-Text "24×94=" "63×33="

Replace-Text "39×57=" "96×37="
Replace-Text "20×52=" "24×99="
Replace-Text "97×45=" "91×36="
Replace-Text "59×60=" "88×90="
Replace-Text "78×17=" "48×50="

Replace-Text "15×68=" "90×50="
Replace-Text "20×75=" "33×35="
Replace-Text "15×45=" "43×71="
Replace-Text "25×35=" "11×12="
Replace-Text "79×32=" "61×31="

Replace-Text "36×41=" "34×63="
Replace-Text "22×55=" "50×52="
Replace-Text "71×67=" "51×65="
Replace-Text "31×47=" "20×49="
Replace-Text "16×21=" "43×65="
